$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M,
# formulas referencing old D/E/etc auto-adjust to the new columns).
$ws.Columns("D:E").Insert()

# --- Row 1 ---
$ws.Range("B1").Value = "bit resolution"
$ws.Range("D1").Value = "clock speed"
$ws.Range("E1").Formula = '="->"'

# --- Row 2 ---
$ws.Range("B2").Value = "day"
$ws.Range("D2").Value = "length of clock tick"
$ws.Range("E2:E4").Formula = '="->"'

# --- Row 4 ---
$ws.Range("B4").Value = "bit resolution"

# --- Row 5 : new "ticks/sec" calculation ---
$ws.Range("F5").Formula = "=1/F4"
$ws.Range("G5").Value = "ticks/sec"

# Column widths (best-fit sizing for the newly inserted columns, closest
# achievable values under this runtime's column-width quantization)
$ws.Range("C1").ColumnWidth = 6.43
$ws.Range("D1").ColumnWidth = 16.8
$ws.Range("E1").ColumnWidth = 1.8

# Restore selection to match the saved workbook state
$ws.Range("F8").Select()
